$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the
#    Heading1 title paragraph ("Play Big Thunder King Strike for
#    Free - Ainsworth Slot Game").
# -----------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range
$metaRange.Text = "Meta description: Explore the immersive jungle theme of Big Thunder King Strike by Ainsworth with free spins, multipliers, and jackpots. Play for free and win big!"

# Bold just the "Meta description" label (first 16 characters).
$boldRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start + 16)
$boldRange.Bold = 1

# Add the leading empty run ("<w:r/>") that matches the style used
# throughout the rest of the document's body paragraphs. Inserting a
# single-paragraph XML fragment consisting of only an empty run at the
# very start of the (already textful) paragraph merges it in cleanly
# as a leading run, without disturbing neighboring paragraphs.
$insertPoint = $d.Range($metaPara.Range.Start, $metaPara.Range.Start)
$xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$emptyRunXml = "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document $xmlNs><w:body><w:p><w:r/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
$insertPoint.InsertXML($emptyRunXml)

# -----------------------------------------------------------------
# 2) Remove the duplicated bold "Play Big Thunder King Strike for
#    Free - Ainsworth Slot Game" paragraph near the end of the
#    document, and update the following italic paragraph's text
#    to the new image-generation prompt.
# -----------------------------------------------------------------
$targetText = "Play Big Thunder King Strike for Free - Ainsworth Slot Game"

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq $targetText -and $i -gt 1) {
        $para.Range.Delete()
        break
    }
}

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newEnd = $lastPara.Range.End - 1
$replaceRange = $d.Range($lastPara.Range.Start, $newEnd)
$replaceRange.Text = 'Please create a feature image for "Big Thunder King Strike" that fits the theme of the game and features a happy Maya warrior wearing glasses. The image should be in a cartoon style and can include elements such as jungle foliage, animals, or tribal weapons. Be creative and use bold colors to capture the game''s eccentric and adventurous vibe.'
